$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right (B11) 5 -> 4, Wrong (C11) -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right (B12) 100 -> 80, Wrong (C12) -2 -> -4, Max text (E12) "100 / 140" -> "76 / 112"
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "76 / 112"
